# Rotate the species-observation rows 2-13 (columns A:AY) according to the
# mapping discovered from the authoritative diff:
#   source row -> destination row (content moves FROM source TO destination)
#   2 -> 4, 3 -> 2, 4 -> 5, 5 -> 6, 6 -> 7, 7 -> 8, 8 -> 9, 9 -> 10,
#   10 -> 11, 11 -> 12, 12 -> 13, 13 -> 3
# This forms a single 12-cycle, so we must buffer all source rows before
# writing any destination row (otherwise we would overwrite data we still
# need to read).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 1   # A
$lastCol  = 51  # AY

$mapping = @{}
$mapping[2]  = 4
$mapping[3]  = 2
$mapping[4]  = 5
$mapping[5]  = 6
$mapping[6]  = 7
$mapping[7]  = 8
$mapping[8]  = 9
$mapping[9]  = 10
$mapping[10] = 11
$mapping[11] = 12
$mapping[12] = 13
$mapping[13] = 3

# Column "I" (Antal / count) is always stored as plain text in this sheet,
# even when its content looks like an integer (e.g. "8", "10"). The engine's
# Value2/Text accessors silently coerce purely-numeric-looking text back to
# a number, so that column needs to be force-typed as text unconditionally.
$alwaysTextCols = @(9)   # I

# --- Phase 1: read every source row (2-13) across columns A:AY into memory,
#     remembering both the value and whether it originated as text (so we
#     can re-apply "Text" formatting and avoid values such as "2019-10-09"
#     or "8" being reinterpreted as dates/numbers on write-back). ---
$bufferValue = @{}
$bufferIsText = @{}
foreach ($srcRow in $mapping.Keys) {
    $rowValues = @{}
    $rowIsText = @{}
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell = $ws.Cells.Item($srcRow, $col)
        $v = $cell.Value2
        $rowValues[$col] = $v
        $isText = ($v -is [string]) -and ($v -ne "")
        if ($alwaysTextCols -contains $col -and $v -ne $null -and "$v" -ne "") {
            $isText = $true
        }
        $rowIsText[$col] = $isText
    }
    $bufferValue[$srcRow] = $rowValues
    $bufferIsText[$srcRow] = $rowIsText
}

# --- Phase 2: write the buffered rows into their destination rows ---
foreach ($srcRow in $mapping.Keys) {
    $dstRow = $mapping[$srcRow]
    $rowValues = $bufferValue[$srcRow]
    $rowIsText = $bufferIsText[$srcRow]
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $destCell = $ws.Cells.Item($dstRow, $col)
        $val = $rowValues[$col]

        if ($rowIsText[$col]) {
            # Force text so numeric-looking or date-looking strings (e.g.
            # "8", "2019-10-09") are preserved as plain text, matching the
            # original inlineStr typing instead of being reinterpreted.
            $destCell.NumberFormat = "@"
            $destCell.Value = "$val"
        } else {
            $destCell.Value = $val
        }
    }
}
